$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row shared strings: "_old" -> "_FV2310", "_new" -> "_FV2404" ---
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value2 = ($cell.Value2 -replace "_old$", "_FV2310")
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value2 = ($cell.Value2 -replace "_new$", "_FV2404")
}

# --- Turn the data range into an Excel Table (ListObject) ---
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U75"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# --- Freeze the header row (row 1) ---
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
